$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, matching style/formatting of existing header row (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-10
$saveValues = @(1, 1, 1, 0, 0, 0, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
